# "Test Commit From master"
# Sheet1 gains a new row: "TestSheetOneFromMaster" is inserted above the
# existing "TestTwo" row, pushing it from A2 down to A3. Column A is then
# auto-sized to fit the new (longer) text, and the saved selection moves
# to B6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 2 (shifts the existing "TestTwo" row down to row 3)
$ws.Range("A2").EntireRow.Insert()

# Populate the newly inserted cell with the new shared string
$ws.Range("A2").Value = "TestSheetOneFromMaster"

# Column A auto-fits to the widest value now in the column
$ws.Columns("A").AutoFit()

# Record the sheet's active selection as saved in the file
$ws.Range("B6").Select()
